# Apply the "add slides about Parser" workbook edit:
#   - re-purpose the sheet: row4 ("gas loss factor") content moves to row5,
#     row4 becomes a new "primary product" / "Oil" row (right-aligned value)
#   - remove the old styled "Production" header row (row6) completely
#   - add a new "Fee" section (rows 12-14: Fee / exploration / post exploration)
#   - add new named ranges (ExplorationFee, PostExplorationFee, PrimaryProduct,
#     Years) and repoint GasLossFactor at its new location
#   - move the active selection to C13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- push the existing "gas loss factor" row down from row 4 to row 5 ---
$ws.Range("A5").Value = $ws.Range("A4").Value()
$ws.Range("B5").Value = $ws.Range("B4").Value()

# --- remove the old bold/shaded "Production" header row (row 6) entirely ---
$ws.Range("A6:J6").Clear()
$ws.Rows("6:6").AutoFit()

# --- row 4 now holds the new "primary product" / "Oil" pair ---
$ws.Range("A4").Value = "primary product"
$ws.Range("B4").Value = "Oil"
$ws.Range("B4").HorizontalAlignment = -4152   # xlRight

# --- new Fee section ---
$ws.Range("A12").Value = "Fee"
$ws.Range("A13").Value = "exploration"
$ws.Range("B13").Value = 1.4
$ws.Range("A14").Value = "post exploration"
$ws.Range("B14").Value = 5.8

# --- defined names ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "GasLossFactor") {
        $n.RefersTo = '=Sheet1!$B$5'
    }
}
$wb.Names.Add("PrimaryProduct", '=Sheet1!$B$4')
$wb.Names.Add("Years", '=Sheet1!$B$2:$J$2')
$wb.Names.Add("ExplorationFee", '=Sheet1!$B$13')
$wb.Names.Add("PostExplorationFee", '=Sheet1!$B$14')

# --- move selection to C13 ---
$ws.Range("C13").Select()
